$wb = $excel.ActiveWorkbook

# --- Fix up selection on the existing Fulfilment_view sheet (P7 -> P19) ---
$wsFulfil = $wb.Worksheets.Item("Fulfilment_view")
$wsFulfil.Range("P19").Select() | Out-Null

# --- Add the new "Thought_leadership_demonstratio" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Thought_leadership_demonstratio"

# Column widths (ColumnWidth is stored as ~ xml_width - 5/6)
$ws4.Columns.Item(5).ColumnWidth = 38 - 5/6
$ws4.Columns.Item(6).ColumnWidth = 8.796875 - 5/6
$ws4.Columns.Item(9).ColumnWidth = 30.6640625 - 5/6
$ws4.Columns.Item(12).ColumnWidth = 22.265625 - 5/6

# Header row
$ws4.Range("C6").Value = "id"
$ws4.Range("E6").Value = "label"
$ws4.Range("G6").Value = "percentage"
$ws4.Range("I6").Value = "title"
$ws4.Range("L6").Value = "heading"

# Row 8
$ws4.Range("C8").Value = 1
$ws4.Range("E8").Value = "TOTAL DEMANDS"
$ws4.Range("G8").Value = 45
$ws4.Range("L8").Value = "Total Revenue: 2023-2024"

# Row 9
$ws4.Range("C9").Value = 2
$ws4.Range("E9").Value = "Intial Net Billed Resources"
$ws4.Range("G9").Value = 45
$ws4.Range("I9").Value = "Netbilled HC UPTICK"
$ws4.Range("L9").Value = "Total Revenue: 2023-2024"

# Row 10
$ws4.Range("C10").Value = 3
$ws4.Range("E10").Value = "Total No Of Demands Filled By Internal / External"
$ws4.Range("I10").Value = "Fulfillment Chanel Performance"

# Row 11
$ws4.Range("C11").Value = 4
$ws4.Range("E11").Value = "No. Of External Fulfilment"
$ws4.Range("G11").Value = 33
$ws4.Range("I11").Value = "Total Demands Vs External Fulfilment"
$ws4.Range("L11").Value = "Total Revenue: 2023-2024"

# Row 12
$ws4.Range("C12").Value = 5
$ws4.Range("E12").Value = "No. Of Internal Fulfilment"
$ws4.Range("G12").Value = 99
$ws4.Range("I12").Value = "Total Demands Vs Internal Fulfilment"
$ws4.Range("L12").Value = "Total Revenue: 2023-2024"

# Selection / active cell on the new sheet, and make it the active (last-viewed) sheet/tab
$ws4.Range("G12").Select() | Out-Null
